# Auto-generated edit script: refresh market-price columns (H-N) on all Leve profit sheets.
# Source: scheduled runner data refresh (chore: update Sheets via scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3924.75
$ws.Range("J43").Value = 4066.3333
$ws.Range("L43").Value = 4066.3333
$ws.Range("N43").Value = -4204.3333
$ws.Range("H64").Value = 8999.799999999999
$ws.Range("I64").Value = 8999.799999999999
$ws.Range("K64").Value = 8999.799999999999
$ws.Range("M64").Value = -8751.799999999999
$ws.Range("H67").Value = 8999.799999999999
$ws.Range("I67").Value = 8999.799999999999
$ws.Range("K67").Value = 8999.799999999999
$ws.Range("M67").Value = -8141.799999999999
$ws.Range("H132").Value = 1223.6833
$ws.Range("I132").Value = 921.12067
$ws.Range("K132").Value = 2763.36201
$ws.Range("M132").Value = -233.3620099999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 18002418
$ws.Range("I45").Value = 3481
$ws.Range("K45").Value = 3481
$ws.Range("M45").Value = -3104
$ws.Range("H49").Value = 40000
$ws.Range("J49").Value = 40000
$ws.Range("L49").Value = 40000
$ws.Range("N49").Value = -40520
$ws.Range("H109").Value = 40188
$ws.Range("I109").Value = 49999
$ws.Range("K109").Value = 49999
$ws.Range("M109").Value = -48612
$ws.Range("H110").Value = 1530
$ws.Range("I110").Value = 1200
$ws.Range("K110").Value = 1200
$ws.Range("M110").Value = 845
$ws.Range("H123").Value = 47214.5
$ws.Range("J123").Value = 47214.5
$ws.Range("L123").Value = 47214.5
$ws.Range("N123").Value = -57014.5
$ws.Range("H128").Value = 70000
$ws.Range("J128").Value = 70000
$ws.Range("L128").Value = 70000
$ws.Range("N128").Value = -79960
$ws.Range("H130").Value = 19826.666
$ws.Range("J130").Value = 19826.666
$ws.Range("L130").Value = 19826.666
$ws.Range("N130").Value = -29866.666
$ws.Range("H138").Value = 69999
$ws.Range("J138").Value = 69999
$ws.Range("L138").Value = 69999
$ws.Range("N138").Value = -80279
$ws.Range("H140").Value = 48998.5
$ws.Range("J140").Value = 48998.5
$ws.Range("L140").Value = 48998.5
$ws.Range("N140").Value = -59358.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3207.5652
$ws.Range("J31").Value = 4991.625
$ws.Range("L31").Value = 4991.625
$ws.Range("N31").Value = -5581.625
$ws.Range("H34").Value = 3207.5652
$ws.Range("J34").Value = 4991.625
$ws.Range("L34").Value = 4991.625
$ws.Range("N34").Value = -5395.625
$ws.Range("H64").Value = 53000
$ws.Range("J64").Value = 53000
$ws.Range("L64").Value = 53000
$ws.Range("N64").Value = -53496
$ws.Range("H67").Value = 53000
$ws.Range("J67").Value = 53000
$ws.Range("L67").Value = 53000
$ws.Range("N67").Value = -54716
$ws.Range("H99").Value = 22224764
$ws.Range("I99").Value = 27780142
$ws.Range("J99").Value = 3250
$ws.Range("K99").Value = 27780142
$ws.Range("L99").Value = 3250
$ws.Range("M99").Value = -27778644
$ws.Range("N99").Value = -6246
$ws.Range("H108").Value = 78684
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 78684
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 78684
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -86364
$ws.Range("H126").Value = 22224764
$ws.Range("I126").Value = 27780142
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 83340426
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -83337956
$ws.Range("N126").Value = -14690
$ws.Range("H132").Value = 1978681.8
$ws.Range("I132").Value = 1978681.8
$ws.Range("K132").Value = 5936045.4
$ws.Range("M132").Value = -5933515.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 371.83334
$ws.Range("I17").Value = 288
$ws.Range("J17").Value = 539.5
$ws.Range("K17").Value = 864
$ws.Range("L17").Value = 1618.5
$ws.Range("M17").Value = -695
$ws.Range("N17").Value = -1956.5
$ws.Range("H56").Value = 6308.143
$ws.Range("I56").Value = 6308.143
$ws.Range("K56").Value = 6308.143
$ws.Range("M56").Value = -5778.143
$ws.Range("H136").Value = 1284.6923
$ws.Range("I136").Value = 1284.6923
$ws.Range("K136").Value = 3854.0769
$ws.Range("M136").Value = 1245.9231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 88887
$ws.Range("J64").Value = 88887
$ws.Range("L64").Value = 88887
$ws.Range("N64").Value = -89383
$ws.Range("H67").Value = 88887
$ws.Range("J67").Value = 88887
$ws.Range("L67").Value = 88887
$ws.Range("N67").Value = -90603
$ws.Range("H97").Value = 1084.8462
$ws.Range("I97").Value = 508.83334
$ws.Range("K97").Value = 508.83334
$ws.Range("M97").Value = -12.83334000000002
$ws.Range("H102").Value = 2627.8
$ws.Range("I102").Value = 2475.3333
$ws.Range("K102").Value = 2475.3333
$ws.Range("M102").Value = -853.3332999999998
$ws.Range("H122").Value = 6497.8335
$ws.Range("I122").Value = 8282.286
$ws.Range("J122").Value = 3999.6
$ws.Range("K122").Value = 24846.858
$ws.Range("L122").Value = 11998.8
$ws.Range("M122").Value = -22396.858
$ws.Range("N122").Value = -16898.8
$ws.Range("H124").Value = 98000
$ws.Range("J124").Value = 98000
$ws.Range("L124").Value = 98000
$ws.Range("N124").Value = -107820
$ws.Range("H126").Value = 4683.647
$ws.Range("I126").Value = 2214.1428
$ws.Range("J126").Value = 6412.3
$ws.Range("K126").Value = 6642.428400000001
$ws.Range("L126").Value = 19236.9
$ws.Range("M126").Value = -4172.428400000001
$ws.Range("N126").Value = -24176.9
$ws.Range("H132").Value = 7423.222
$ws.Range("I132").Value = 8519.666999999999
$ws.Range("K132").Value = 25559.001
$ws.Range("M132").Value = -23029.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6042
$ws.Range("I22").Value = 736.1818
$ws.Range("J22").Value = 13337.5
$ws.Range("K22").Value = 736.1818
$ws.Range("L22").Value = 13337.5
$ws.Range("M22").Value = -441.1818
$ws.Range("N22").Value = -13927.5
$ws.Range("H27").Value = 6042
$ws.Range("I27").Value = 736.1818
$ws.Range("J27").Value = 13337.5
$ws.Range("K27").Value = 736.1818
$ws.Range("L27").Value = 13337.5
$ws.Range("M27").Value = -629.1818
$ws.Range("N27").Value = -13551.5
$ws.Range("H40").Value = 22225352
$ws.Range("I40").Value = 5301.3335
$ws.Range("K40").Value = 5301.3335
$ws.Range("M40").Value = -5165.3335
$ws.Range("H68").Value = 9001.5
$ws.Range("I68").Value = 6000.5
$ws.Range("J68").Value = 12002.5
$ws.Range("K68").Value = 6000.5
$ws.Range("L68").Value = 12002.5
$ws.Range("M68").Value = -5251.5
$ws.Range("N68").Value = -13500.5
$ws.Range("H71").Value = 9001.5
$ws.Range("I71").Value = 6000.5
$ws.Range("J71").Value = 12002.5
$ws.Range("K71").Value = 30002.5
$ws.Range("L71").Value = 60012.5
$ws.Range("M71").Value = -26258.5
$ws.Range("N71").Value = -67500.5
$ws.Range("H122").Value = 70591880
$ws.Range("I122").Value = 76926480
$ws.Range("K122").Value = 230779440
$ws.Range("M122").Value = -230776990
$ws.Range("H129").Value = 73064
$ws.Range("J129").Value = 72398.8
$ws.Range("L129").Value = 72398.8
$ws.Range("N129").Value = -82398.8
$ws.Range("H132").Value = 1827.0667
$ws.Range("I132").Value = 1723.9231
$ws.Range("K132").Value = 5171.7693
$ws.Range("M132").Value = -2641.7693
$ws.Range("H133").Value = 51454.668
$ws.Range("J133").Value = 55745.6
$ws.Range("L133").Value = 55745.6
$ws.Range("N133").Value = -60805.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 100000
$ws.Range("J16").Value = 100000
$ws.Range("L16").Value = 100000
$ws.Range("N16").Value = -100584
$ws.Range("H62").Value = 4300
$ws.Range("I62").Value = 4066.6667
$ws.Range("K62").Value = 4066.6667
$ws.Range("M62").Value = -3442.6667
$ws.Range("H65").Value = 4300
$ws.Range("I65").Value = 4066.6667
$ws.Range("K65").Value = 20333.3335
$ws.Range("M65").Value = -17213.3335
$ws.Range("H86").Value = 53333
$ws.Range("J86").Value = 53333
$ws.Range("L86").Value = 53333
$ws.Range("N86").Value = -55579
$ws.Range("H89").Value = 53333
$ws.Range("J89").Value = 53333
$ws.Range("L89").Value = 266665
$ws.Range("N89").Value = -277897
$ws.Range("H107").Value = 4700.4
$ws.Range("I107").Value = 2833
$ws.Range("K107").Value = 8499
$ws.Range("M107").Value = -6579
$ws.Range("H109").Value = 61111
$ws.Range("J109").Value = 61111
$ws.Range("L109").Value = 61111
$ws.Range("N109").Value = -63885
$ws.Range("H122").Value = 2439.7727
$ws.Range("I122").Value = 1618
$ws.Range("K122").Value = 4854
$ws.Range("M122").Value = -2404
$ws.Range("H136").Value = 1644.9375
$ws.Range("I136").Value = 1644.9375
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4934.8125
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2384.8125
$ws.Range("N136").ClearContents()
